# edit.ps1 -- PowerShell-style PowerPoint COM-interop script
#
# Reproduces the commit:
#   1) The table on slide 16 is switched to a different (built-in) table
#      style - GUID {A10F5928-0615-4EDB-A427-199EC96E0A85} (the package's
#      default "Table_0" style) becomes {374B2A09-7F2C-4E4A-8F76-C3AE0ED87230}.
#   2) The deck's theme ("Integral") is recoloured to the stock PowerPoint
#      "Office" palette (dk2/lt2/accent1-6/hlink/folHlink all change; dk1
#      and lt1 - pure black/white - are identical in both palettes already).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style on slide 16 (3rd shape on that slide is the table).
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{374B2A09-7F2C-4E4A-8F76-C3AE0ED87230}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the theme from "Integral" to the stock "Office" palette.
# ---------------------------------------------------------------------------
function ConvertTo-BgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# ThemeColorScheme item order is: dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink.
$officePalette = @{
    3  = "44546A"   # dk2
    4  = "E7E6E6"   # lt2
    5  = "5B9BD5"   # accent1
    6  = "ED7D31"   # accent2
    7  = "A5A5A5"   # accent3
    8  = "FFC000"   # accent4
    9  = "4472C4"   # accent5
    10 = "70AD47"   # accent6
    11 = "0563C1"   # hlink
    12 = "954F72"   # folHlink
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
foreach ($idx in $officePalette.Keys) {
    $themeColors.Item($idx).RGB = ConvertTo-BgrInt $officePalette[$idx]
}
